# Fix the username for the newly-created Priya Giri "Business" account:
# it was entered as the raw login "priya.giri" instead of following the
# established "Firstname.Lastname.Role" naming convention used by every
# other row (e.g. "Andrew.Nisbet.Business").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A20").Value = "Priya.Giri.Business"

# Leave the selection where the author ended up after fixing the row.
$ws.Range("A20").Select()
